$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 139, shifting rows 139:186 down to 140:187
$ws.Rows.Item(139).Insert()

# Populate the newly inserted row 139 with the new record's data
$ws.Range("A139").Value = 3
$ws.Range("B139").Value = "Femacal de La Calera"
$ws.Range("C139").Value = "Coquimbo"
$ws.Range("D139").Value = 44588
$ws.Range("E139").Value = 5
$ws.Range("F139").Value = 100112010
$ws.Range("G139").Value = "Achicoria"
$ws.Range("H139").Value = "Sin especificar"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 60
$ws.Range("K139").Value = 6000
$ws.Range("L139").Value = 6000
$ws.Range("M139").Value = 6000
$ws.Range("N139").Value = "`$/caja 16 unidades"
$ws.Range("O139").Value = "Provincia de Quillota"
$ws.Range("P139").Value = 375
$ws.Range("Q139").Value = 16
$ws.Range("R139").Value = "Hortaliza"
